# Update the "Förändrad" (Changed) date in column C for data rows 2-21
# from 45170 (2023-09-01) to 45174 (2023-09-05).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45170) {
        $cell.Value = 45174
    }
}
